# docs : Plan week 3
# Adds a new "Week 4" worksheet (copied/adapted from "Week 3") with the
# week-4 task assignments, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create "Week 4" as a copy of "Week 3" (placed right after it) so
#    it inherits the same column widths / merged layout / styles.
# ------------------------------------------------------------------
$weekCount = $wb.Worksheets.Count
$ws3 = $wb.Worksheets.Item($weekCount)
$ws3.Copy($null, $ws3)
$ws4 = $wb.Worksheets.Item($weekCount + 1)
$ws4.Name = "Week 4"

# ------------------------------------------------------------------
# 2. Block 1 (Hoàng Thị Thảo Nhi) - rows 2:4
# ------------------------------------------------------------------
$ws4.Range("B2").Value = "Thiết kế front end giao diện đăng nhập, đăng kí"
$ws4.Range("C2").Value = 45667
$ws4.Range("C3").HorizontalAlignment = 1

# ------------------------------------------------------------------
# 3. Block 2 (Nguyễn Hoài Nam) - rows 5:7
# ------------------------------------------------------------------
$ws4.Range("B5").Value = "Thiết kế front end cho trang sản phẩm"
$ws4.Range("C5").Value = 45667
$ws4.Range("B6").ClearContents()
$ws4.Range("C6").HorizontalAlignment = 1

# ------------------------------------------------------------------
# 4. Block 3 (Đặng Đức Huy) - rows 8:10
# ------------------------------------------------------------------
$ws4.Range("B8").Value = "Làm chức năng Tìm kiếm"
$ws4.Range("C8").Value = 45667
$ws4.Range("B9").Value = "Nạp data vào website để có thể hiển thị sản phẩm"
$ws4.Range("B10").ClearContents()

# C8:C10 was one merged date cell in "Week 3"; in "Week 4" only C8:C9
# are merged and C10 stands alone (blank, vertical-centered only).
$ws4.Range("C8:C10").UnMerge()
$ws4.Range("C8:C9").Merge()
$ws4.Range("C10").ClearContents()
$ws4.Range("C10").HorizontalAlignment = 1

# ------------------------------------------------------------------
# 5. Selection bookkeeping + make "Week 4" the active/visible tab.
# ------------------------------------------------------------------
$ws3.Range("B2").Select()
$ws4.Range("B16").Select()
$ws4.Activate()
